# Daily status update - 2/21 - Yash
# Yash (row 8) reports availability of 50% for tomorrow, matching the
# percentage formatting already used for Srinandan's entry in C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give C8 the same percentage number format as C7, then set its value to 50%.
$ws.Range("C8").NumberFormat = $ws.Range("C7").NumberFormat
$ws.Range("C8").Value = 0.5

# Move the active selection to C8 (where the update was made).
$ws.Range("C8").Select()
